$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1792.6578
$ws.Range("I40").Value = 1823.1786
$ws.Range("J40").Value = 1707.2
$ws.Range("K40").Value = 1823.1786
$ws.Range("L40").Value = 1707.2
$ws.Range("M40").Value = -1648.1786
$ws.Range("N40").Value = -2057.2

$ws.Range("H64").Value = 3910
$ws.Range("I64").Value = 3887.5
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3887.5
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3639.5
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 3910
$ws.Range("I67").Value = 3887.5
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3887.5
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -3029.5
$ws.Range("N67").Value = -5716

$ws.Range("H116").Value = 6038.231
$ws.Range("I116").Value = 9277.615
$ws.Range("J116").Value = 2798.8462
$ws.Range("K116").Value = 9277.615
$ws.Range("L116").Value = 2798.8462
$ws.Range("M116").Value = -5835.615
$ws.Range("N116").Value = -9682.8462

$ws.Range("H128").Value = 21432.223
$ws.Range("J128").Value = 21432.223
$ws.Range("L128").Value = 21432.223
$ws.Range("N128").Value = -31392.223

$ws.Range("H138").Value = 3512.5168
$ws.Range("I138").Value = 805.60785
$ws.Range("J138").Value = 7145.4736
$ws.Range("K138").Value = 2416.82355
$ws.Range("L138").Value = 21436.4208
$ws.Range("M138").Value = 2723.17645
$ws.Range("N138").Value = -31716.4208

$ws.Range("H139").Value = 29166.666
$ws.Range("J139").Value = 29166.666
$ws.Range("L139").Value = 29166.666
$ws.Range("N139").Value = -39446.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("N60").Value = 0

$ws.Range("H61").Value = 3364.34
$ws.Range("I61").Value = 3296.1875
$ws.Range("K61").Value = 3296.1875
$ws.Range("M61").Value = -3084.1875

$ws.Range("H74").Value = 13159076
$ws.Range("I74").Value = 1077.4517
$ws.Range("J74").Value = 71430216
$ws.Range("K74").Value = 1077.4517
$ws.Range("L74").Value = 71430216
$ws.Range("M74").Value = -203.4517000000001
$ws.Range("N74").Value = -71431964

$ws.Range("H77").Value = 13159076
$ws.Range("I77").Value = 1077.4517
$ws.Range("J77").Value = 71430216
$ws.Range("K77").Value = 5387.2585
$ws.Range("L77").Value = 357151080
$ws.Range("M77").Value = -1019.2585
$ws.Range("N77").Value = -357159816

$ws.Range("H110").Value = 602.8570999999999
$ws.Range("J110").Value = 210
$ws.Range("L110").Value = 210
$ws.Range("N110").Value = -4300

$ws.Range("H132").Value = 1967.1094
$ws.Range("I132").Value = 1140.3062
$ws.Range("K132").Value = 3420.9186
$ws.Range("M132").Value = -890.9186

$ws.Range("H133").Value = 24800
$ws.Range("J133").Value = 24800
$ws.Range("L133").Value = 24800
$ws.Range("N133").Value = -29860

$ws.Range("H136").Value = 3364.34
$ws.Range("I136").Value = 3296.1875
$ws.Range("K136").Value = 9888.5625
$ws.Range("M136").Value = -7338.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1400.2
$ws.Range("I107").Value = 1497
$ws.Range("K107").Value = 1497
$ws.Range("M107").Value = 423

$ws.Range("H134").Value = 3464
$ws.Range("I134").Value = 3583.7021
$ws.Range("J134").Value = 2901.4
$ws.Range("K134").Value = 10751.1063
$ws.Range("L134").Value = 8704.200000000001
$ws.Range("M134").Value = -8216.106299999999
$ws.Range("N134").Value = -13774.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8923.076999999999
$ws.Range("I4").Value = 49000
$ws.Range("J4").Value = 7320
$ws.Range("K4").Value = 49000
$ws.Range("L4").Value = 7320
$ws.Range("M4").Value = -48888
$ws.Range("N4").Value = -7544

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7143806.5
$ws.Range("I131").Value = 25000532
$ws.Range("J131").Value = 1116.4
$ws.Range("K131").Value = 75001596
$ws.Range("L131").Value = 3349.2
$ws.Range("M131").Value = -74996556
$ws.Range("N131").Value = -13429.2

$ws.Range("H132").Value = 7408707
$ws.Range("I132").Value = 1018.1818
$ws.Range("J132").Value = 27779850
$ws.Range("K132").Value = 9163.636199999999
$ws.Range("L132").Value = 250018650
$ws.Range("M132").Value = -6633.636199999999
$ws.Range("N132").Value = -250023710

$ws.Range("H139").Value = 5401.6216
$ws.Range("I139").Value = 11632
$ws.Range("J139").Value = 3094.074
$ws.Range("K139").Value = 34896
$ws.Range("L139").Value = 9282.222
$ws.Range("M139").Value = -29756
$ws.Range("N139").Value = -19562.222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6922.4
$ws.Range("I126").Value = 8996.214
$ws.Range("J126").Value = 2083.5
$ws.Range("K126").Value = 26988.642
$ws.Range("L126").Value = 6250.5
$ws.Range("M126").Value = -24518.642
$ws.Range("N126").Value = -11190.5

$ws.Range("H132").Value = 4275714
$ws.Range("I132").Value = 5210182.5
$ws.Range("J132").Value = 3856.8572
$ws.Range("K132").Value = 15630547.5
$ws.Range("L132").Value = 11570.5716
$ws.Range("M132").Value = -15628017.5
$ws.Range("N132").Value = -16630.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 850
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -630
$ws.Range("N16").Value = -1240

$ws.Range("H46").Value = 22222950
$ws.Range("I46").Value = 166667150
$ws.Range("J46").Value = 766.8461
$ws.Range("K46").Value = 166667150
$ws.Range("L46").Value = 766.8461
$ws.Range("M46").Value = -166666962
$ws.Range("N46").Value = -1142.8461

$ws.Range("H93").Value = 21739862
$ws.Range("I93").Value = 835.3333
$ws.Range("J93").Value = 62500536
$ws.Range("K93").Value = 835.3333
$ws.Range("L93").Value = 62500536
$ws.Range("M93").Value = 412.6667
$ws.Range("N93").Value = -62503032

$ws.Range("H132").Value = 8188514
$ws.Range("I132").Value = 11462830
$ws.Range("J132").Value = 2725.0625
$ws.Range("K132").Value = 34388490
$ws.Range("L132").Value = 8175.1875
$ws.Range("M132").Value = -34385960
$ws.Range("N132").Value = -13235.1875

$ws.Range("H136").Value = 9958.75
$ws.Range("I136").Value = 9035.352999999999
$ws.Range("J136").Value = 12201.286
$ws.Range("K136").Value = 27106.059
$ws.Range("L136").Value = 36603.858
$ws.Range("M136").Value = -24556.059
$ws.Range("N136").Value = -41703.858

$ws.Range("H140").Value = 50756.816
$ws.Range("J140").Value = 50756.816
$ws.Range("L140").Value = 50756.816
$ws.Range("N140").Value = -61116.816

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 333333470
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 1006.5
$ws.Range("I132").Value = 712.0417
$ws.Range("K132").Value = 2136.1251
$ws.Range("M132").Value = 393.8748999999998

$ws.Range("H133").Value = 45357.5
$ws.Range("J133").Value = 45357.5
$ws.Range("L133").Value = 45357.5
$ws.Range("N133").Value = -55477.5

$ws.Range("H136").Value = 3877125.8
$ws.Range("I136").Value = 703.4211
$ws.Range("J136").Value = 6945960
$ws.Range("K136").Value = 2110.2633
$ws.Range("L136").Value = 20837880
$ws.Range("M136").Value = 439.7366999999999
$ws.Range("N136").Value = -20842980

$ws.Range("H140").Value = 37543
$ws.Range("J140").Value = 37543
$ws.Range("L140").Value = 37543
$ws.Range("N140").Value = -47903
